# Apply content edits to the "about_app" translations sheet.
# Only three cell values actually change text content between the
# before/after OOXML (everything else is shared-string-table bookkeeping
# that Excel recomputes automatically when the workbook is saved):
#   C23 -> updated Dutch wording for the "send my APP ID" question
#   B46 -> English "about organizations" paragraph gains a <a ...>here</a> link
#   C46 -> Dutch "about organizations" paragraph gains a <a ...>hier</a> link
#          (and its stray non-breaking spaces are normalised to regular spaces)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("about_app")

$ws.Range("C23").Value = "Ik heb zojuist een e-mail ontvangen van de Universiteit Leiden met de vraag om mijn APP-ID toe te sturen, hoe doe ik dat?"

$ws.Range("B46").Value = 'For this study, the research team at Leiden University is using this App. This App is operated by the University of Zurich Spin-off QuantActions Ltd (Lausanne, Switzerland) and you can find its detailed privacy policy <a target=_blank_  href="https://quantactions.com/privacy/">here</a>. Google PlayStore is used to distribute this App to users like you.'

$ws.Range("C46").Value = 'Het onderzoeksteam van de Universiteit Leiden gebruikt deze app voor dit onderzoek. Deze app wordt beheerd door de University of Zurich Spin-off QuantActions Ltd (Lausanne, Zwitserland) en u kunt het gedetailleerde privacybeleid <a target=_blank_  href="https://quantactions.com/privacy/">hier</a> vinden. Google PlayStore wordt gebruikt om deze app te distribueren naar gebruikers zoals u. '

# Reflect where the author last clicked before saving.
$ws.Range("C49").Select()
